$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: "რეგისტრირებული ოჯახი" (registered families) — updated figures for 2015-2021 (columns E:K)
$ws.Range("E4").Value = 6100
$ws.Range("F4").Value = 5026
$ws.Range("G4").Value = 3333
$ws.Range("H4").Value = 3191
$ws.Range("I4").Value = 2913
$ws.Range("J4").Value = 3058
$ws.Range("K4").Value = 3275

# Row 5: "საარსებო შემწეობის მიმღები ოჯახი" (families receiving subsistence allowance) — updated figures for 2015-2021 (columns E:K)
$ws.Range("E5").Value = 696
$ws.Range("F5").Value = 784
$ws.Range("G5").Value = 729
$ws.Range("H5").Value = 685
$ws.Range("I5").Value = 681
$ws.Range("J5").Value = 804
$ws.Range("K5").Value = 1057
